$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# --- Step 1: insert 3 new rows before the current last row (old row 11,
#     "groups: delta CSneu / CS-") so it becomes row 14 and rows 11-13
#     are fresh rows for the new pairwise "imagery"/"classical" comparisons. ---
$beforeRow = $tbl.Rows(11)
$tbl.Rows.Add($beforeRow) | Out-Null
$tbl.Rows.Add($beforeRow) | Out-Null
$tbl.Rows.Add($beforeRow) | Out-Null

# --- Step 2: write final cell text for every data row (rows 3-14). ---

# Row 3: total sample: CS+av vs CS+neu
$tbl.Cell(3,1).Range.Text = "total sample: CS+av vs CS+neu"
$tbl.Cell(3,2).Range.Text = "7.04"
$tbl.Cell(3,3).Range.Text = "47"
$tbl.Cell(3,4).Range.Text = "< .001"
$tbl.Cell(3,5).Range.Text = "1.02"
$tbl.Cell(3,6).Range.Text = "3.51e+06"

# Row 4: total sample: CS+av vs CS-
$tbl.Cell(4,1).Range.Text = "total sample: CS+av vs CS-"
$tbl.Cell(4,2).Range.Text = "6.21"
$tbl.Cell(4,3).Range.Text = "47"
$tbl.Cell(4,4).Range.Text = "< .001"
$tbl.Cell(4,5).Range.Text = "0.90"
$tbl.Cell(4,6).Range.Text = "2.23e+05"

# Row 5: total sample: CSneu vs CS-
$tbl.Cell(5,1).Range.Text = "total sample: CSneu vs CS-"
$tbl.Cell(5,2).Range.Text = "-0.26"
$tbl.Cell(5,3).Range.Text = "47"
$tbl.Cell(5,4).Range.Text = ".796"
$tbl.Cell(5,5).Range.Text = "-0.04"
$tbl.Cell(5,6).Range.Text = "1.62e-01"

# Row 6: between groups: delta CS+av / CS+neu
$tbl.Cell(6,1).Range.Text = "between groups: delta CS+av / CS+neu"
$tbl.Cell(6,2).Range.Text = "2.03"
$tbl.Cell(6,3).Range.Text = "45"
$tbl.Cell(6,4).Range.Text = ".146"
$tbl.Cell(6,5).Range.Text = "0.58"
$tbl.Cell(6,6).Range.Text = "1.48e+00"

# Row 7: between groups: delta CS+av / CS-
$tbl.Cell(7,1).Range.Text = "between groups: delta CS+av / CS-"
$tbl.Cell(7,2).Range.Text = "1.83"
$tbl.Cell(7,3).Range.Text = "45"
$tbl.Cell(7,4).Range.Text = ".221"
$tbl.Cell(7,5).Range.Text = "0.53"
$tbl.Cell(7,6).Range.Text = "1.11e+00"

# Row 8: between groups: delta CSneu / CS-
$tbl.Cell(8,1).Range.Text = "between groups: delta CSneu / CS-"
$tbl.Cell(8,2).Range.Text = "0.00"
$tbl.Cell(8,3).Range.Text = "31"
$tbl.Cell(8,4).Range.Text = "1"
$tbl.Cell(8,5).Range.Text = "0.00"
$tbl.Cell(8,6).Range.Text = "2.87e-01"

# Row 9: imagery: CS+av vs CS+neu
$tbl.Cell(9,1).Range.Text = "imagery: CS+av vs CS+neu"
$tbl.Cell(9,2).Range.Text = "3.50"
$tbl.Cell(9,3).Range.Text = "23"
$tbl.Cell(9,4).Range.Text = "< .001"
$tbl.Cell(9,5).Range.Text = "0.71"
$tbl.Cell(9,6).Range.Text = "3.96e+01"

# Row 10: imagery: CS+av vs CS-
$tbl.Cell(10,1).Range.Text = "imagery: CS+av vs CS-"
$tbl.Cell(10,2).Range.Text = "2.95"
$tbl.Cell(10,3).Range.Text = "23"
$tbl.Cell(10,4).Range.Text = ".004"
$tbl.Cell(10,5).Range.Text = "0.60"
$tbl.Cell(10,6).Range.Text = "1.27e+01"

# Row 11 (new): imagery: CSneu vs CS-
$tbl.Cell(11,1).Range.Text = "imagery: CSneu vs CS-"
$tbl.Cell(11,2).Range.Text = "-0.14"
$tbl.Cell(11,3).Range.Text = "23"
$tbl.Cell(11,4).Range.Text = ".89"
$tbl.Cell(11,5).Range.Text = "-0.03"
$tbl.Cell(11,6).Range.Text = "2.17e-01"

# Row 12 (new): classical: CS+av vs CS+neu
$tbl.Cell(12,1).Range.Text = "classical: CS+av vs CS+neu"
$tbl.Cell(12,2).Range.Text = "7.01"
$tbl.Cell(12,3).Range.Text = "23"
$tbl.Cell(12,4).Range.Text = "< .001"
$tbl.Cell(12,5).Range.Text = "1.43"
$tbl.Cell(12,6).Range.Text = "8.77e+04"

# Row 13 (new): classical: CS+av vs CS-
$tbl.Cell(13,1).Range.Text = "classical: CS+av vs CS-"
$tbl.Cell(13,2).Range.Text = "6.40"
$tbl.Cell(13,3).Range.Text = "23"
$tbl.Cell(13,4).Range.Text = "< .001"
$tbl.Cell(13,5).Range.Text = "1.31"
$tbl.Cell(13,6).Range.Text = "2.34e+04"

# Row 14 (was old row 11): classical: CSneu vs CS-
$tbl.Cell(14,1).Range.Text = "classical: CSneu vs CS-"
$tbl.Cell(14,2).Range.Text = "-0.33"
$tbl.Cell(14,3).Range.Text = "23"
$tbl.Cell(14,4).Range.Text = ".747"
$tbl.Cell(14,5).Range.Text = "-0.07"
$tbl.Cell(14,6).Range.Text = "2.25e-01"

Write-Output "rows now: $($tbl.Rows.Count)"
